$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updated covid_deaths counts for existing rows (dates 2020-12-15 .. 2020-12-22ish) ---
$ws.Range("C997").Value = 40
$ws.Range("C1044").Value = 31
$ws.Range("C1055").Value = 53
$ws.Range("C1061").Value = 7
$ws.Range("C1064").Value = 36
$ws.Range("C1080").Value = 38
$ws.Range("C1085").Value = 49
$ws.Range("C1087").Value = 3
$ws.Range("C1090").Value = 34
$ws.Range("C1095").Value = 14
$ws.Range("C1096").Value = 41
$ws.Range("C1100").Value = 14
$ws.Range("C1101").Value = 48
$ws.Range("C1103").Value = 10
$ws.Range("C1104").Value = 23
$ws.Range("C1105").Value = 35
$ws.Range("C1107").Value = 12
$ws.Range("C1108").Value = 19
$ws.Range("C1109").Value = 36
$ws.Range("C1112").Value = 19
$ws.Range("C1113").Value = 31
$ws.Range("C1115").Value = 15
$ws.Range("C1116").Value = 24
$ws.Range("C1117").Value = 35
$ws.Range("C1120").Value = 2
$ws.Range("C1122").Value = 18
$ws.Range("C1123").Value = 23

# --- Rows 1124-1126 (date 2020-12-23 / serial 44188) had their age-group rows
#     re-ordered / re-labelled and counts changed ---
$ws.Range("B1124").Value = "20-29"

$ws.Range("B1125").Value = "40-49"
$ws.Range("C1125").Value = 3

$ws.Range("B1126").Value = "50-59"
$ws.Range("C1126").Value = 2

# --- New trailing rows 1127-1152: remainder of 2020-12-23 plus 2020-12-24
#     through 2020-12-28 ---
$newRows = @(
    @(1127, 44188, "60-69", 7),
    @(1128, 44188, "70-79", 19),
    @(1129, 44188, "80+",   26),
    @(1130, 44189, "50-59", 4),
    @(1131, 44189, "60-69", 5),
    @(1132, 44189, "70-79", 16),
    @(1133, 44189, "80+",   21),
    @(1134, 44190, "30-39", 1),
    @(1135, 44190, "50-59", 5),
    @(1136, 44190, "60-69", 3),
    @(1137, 44190, "70-79", 6),
    @(1138, 44190, "80+",   32),
    @(1139, 44191, "50-59", 3),
    @(1140, 44191, "60-69", 12),
    @(1141, 44191, "70-79", 12),
    @(1142, 44191, "80+",   25),
    @(1143, 44192, "40-49", 2),
    @(1144, 44192, "50-59", 1),
    @(1145, 44192, "60-69", 7),
    @(1146, 44192, "70-79", 18),
    @(1147, 44192, "80+",   25),
    @(1148, 44193, "40-49", 1),
    @(1149, 44193, "50-59", 4),
    @(1150, 44193, "60-69", 3),
    @(1151, 44193, "70-79", 8),
    @(1152, 44193, "80+",   13)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $dateSerial = $r[1]
    $ageGrp = $r[2]
    $deaths = $r[3]

    $aCell = $ws.Cells.Item($rowNum, 1)
    $aCell.Value = $dateSerial
    $aCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($rowNum, 2).Value = $ageGrp
    $ws.Cells.Item($rowNum, 3).Value = $deaths
}
